$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.178.65"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.033.01"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.43"
$ws.Range("E5").Value = "  -1.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.25"
$ws.Range("E6").Value = "  +3.18%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.029.81"
$ws.Range("E8").Value = "  +0.82%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.520"
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.62"
$ws.Range("E10").Value = "  -0.46%  "
$ws.Range("E11").Value = "  -1.54%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.481"
$ws.Range("E12").Value = "  +5.52%  "
$ws.Range("E13").Value = "  -3.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.79"
$ws.Range("E14").Value = "  +6.41%  "
$ws.Range("E15").Value = "  -0.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.130.82"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.539.37"
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.24"
$ws.Range("E18").Value = "  +4.65%  "
$ws.Range("E19").Value = "  +19.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.035.67"
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "469.84"
$ws.Range("E21").Value = "  +3.33%  "
$ws.Range("E22").Value = "  +2.75%  "
$ws.Range("E23").Value = "  +1.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.26"
$ws.Range("E24").Value = "  +1.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.88"
$ws.Range("E25").Value = "  +4.81%  "
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.04"
$ws.Range("E27").Value = "  -4.66%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.31"
$ws.Range("E29").Value = "  +2.12%  "
$ws.Range("E30").Value = "  +1.84%  "
$ws.Range("E31").Value = "  +0.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.120"
$ws.Range("E32").Value = "  +8.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0₃0991"
$ws.Range("E33").Value = "  -6.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.16"
$ws.Range("E34").Value = "  +3.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.990"
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.84"
$ws.Range("E37").Value = "  +0.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.92"
$ws.Range("E38").Value = "  +8.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.05"
$ws.Range("E39").Value = "  -6.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "49.57"
$ws.Range("E40").Value = "  -0.66%  "
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("E42").Value = "  -1.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.63"
$ws.Range("E43").Value = "  +2.26%  "
$ws.Range("E44").Value = "  -5.92%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0359"
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "380.87"
$ws.Range("E46").Value = "  -3.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.728.97"
$ws.Range("E47").Value = "  -2.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.29"
$ws.Range("E48").Value = "  +0.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.84"
$ws.Range("E50").Value = "  +3.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.23"
$ws.Range("E51").Value = "  +3.64%  "
